$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '26.511.35'
Set-TextValue $ws.Range("E2") '  -0.32%  '
Set-TextValue $ws.Range("D3") '1.721.06'
Set-TextValue $ws.Range("E3") '  -1.34%  '
Set-TextValue $ws.Range("D4") '0.9947'
Set-TextValue $ws.Range("E4") '  -0.52%  '
Set-TextValue $ws.Range("D5") '240.37'
Set-TextValue $ws.Range("E5") '  -2.50%  '
Set-TextValue $ws.Range("D6") '0.9956'
Set-TextValue $ws.Range("E6") '  -0.46%  '
Set-TextValue $ws.Range("D7") '0.4903'
Set-TextValue $ws.Range("E7") '  -0.38%  '
Set-TextValue $ws.Range("E8") '  -3.32%  '
Set-TextValue $ws.Range("D9") '0.06188'
Set-TextValue $ws.Range("E9") '  -1.59%  '
Set-TextValue $ws.Range("D10") '1.724.06'
Set-TextValue $ws.Range("E10") '  -1.09%  '
Set-TextValue $ws.Range("D11") '0.06948'
Set-TextValue $ws.Range("E11") '  -1.37%  '
Set-TextValue $ws.Range("D12") '15.60'
Set-TextValue $ws.Range("E12") '  -0.93%  '
Set-TextValue $ws.Range("D13") '0.6043'
Set-TextValue $ws.Range("E13") '  -1.85%  '
Set-TextValue $ws.Range("D14") '4.463'
Set-TextValue $ws.Range("E14") '  -2.66%  '
Set-TextValue $ws.Range("D15") '76.63'
Set-TextValue $ws.Range("E15") '  -1.94%  '
Set-TextValue $ws.Range("D16") '0.9959'
Set-TextValue $ws.Range("E16") '  -0.45%  '
Set-TextValue $ws.Range("D17") '26.336.35'
Set-TextValue $ws.Range("D18") '0.9944'
Set-TextValue $ws.Range("E18") '  -0.59%  '
Set-TextValue $ws.Range("D19") '0.000007116'
Set-TextValue $ws.Range("E19") '  -2.65%  '
Set-TextValue $ws.Range("D20") '11.30'
Set-TextValue $ws.Range("E20") '  -2.31%  '
Set-TextValue $ws.Range("D21") '1.941.15'
Set-TextValue $ws.Range("E21") '  -1.34%  '
Set-TextValue $ws.Range("D22") '4.403'
Set-TextValue $ws.Range("E22") '  -3.79%  '
Set-TextValue $ws.Range("D23") '8.420'
Set-TextValue $ws.Range("E23") '  -3.47%  '
Set-TextValue $ws.Range("D24") '5.074'
Set-TextValue $ws.Range("E24") '  -3.65%  '
Set-TextValue $ws.Range("D25") '137.91'
Set-TextValue $ws.Range("E25") '  -1.22%  '
Set-TextValue $ws.Range("D26") '15.22'
Set-TextValue $ws.Range("E26") '  -1.58%  '
Set-TextValue $ws.Range("E27") '  -2.62%  '
Set-TextValue $ws.Range("D28") '1.743'
Set-TextValue $ws.Range("E28") '  -1.18%  '
Set-TextValue $ws.Range("D29") '105.64'
Set-TextValue $ws.Range("E29") '  -1.96%  '
Set-TextValue $ws.Range("D30") '3.906'
Set-TextValue $ws.Range("E30") '  -3.56%  '
Set-TextValue $ws.Range("D31") '0.07919'
Set-TextValue $ws.Range("E31") '  -1.46%  '
Set-TextValue $ws.Range("D32") '3.624'
Set-TextValue $ws.Range("E32") '  -3.16%  '
Set-TextValue $ws.Range("D33") '0.04475'
Set-TextValue $ws.Range("E33") '  -3.31%  '
Set-TextValue $ws.Range("D34") '2.603'
Set-TextValue $ws.Range("E34") '  -0.35%  '
Set-TextValue $ws.Range("D35") '0.9973'
Set-TextValue $ws.Range("E35") '  -2.10%  '
Set-TextValue $ws.Range("D36") '0.6166'
Set-TextValue $ws.Range("E36") '  -3.42%  '
Set-TextValue $ws.Range("D37") '0.9481'
Set-TextValue $ws.Range("E37") '  +5.52%  '
Set-TextValue $ws.Range("D38") '2.002'
Set-TextValue $ws.Range("E38") '  -3.39%  '
Set-TextValue $ws.Range("D39") '2.391'
Set-TextValue $ws.Range("E39") '  -1.24%  '
Set-TextValue $ws.Range("D40") '0.9950'
Set-TextValue $ws.Range("E40") '  -0.83%  '
Set-TextValue $ws.Range("D41") '0.01486'
Set-TextValue $ws.Range("E41") '  -1.21%  '
Set-TextValue $ws.Range("D42") '99.64'
Set-TextValue $ws.Range("E42") '  -2.26%  '
Set-TextValue $ws.Range("D43") '5.450'
Set-TextValue $ws.Range("E43") '  +0.37%  '
Set-TextValue $ws.Range("D44") '0.3817'
Set-TextValue $ws.Range("E44") '  -2.69%  '
Set-TextValue $ws.Range("D45") '6.890'
Set-TextValue $ws.Range("E45") '  +0.26%  '
Set-TextValue $ws.Range("D46") '0.1156'
Set-TextValue $ws.Range("E46") '  -2.22%  '
Set-TextValue $ws.Range("D47") '0.05358'
Set-TextValue $ws.Range("E47") '  -0.72%  '
Set-TextValue $ws.Range("E48") '  -0.36%  '
Set-TextValue $ws.Range("D49") '7.709'
Set-TextValue $ws.Range("E49") '  -0.85%  '
Set-TextValue $ws.Range("D50") '51.25'
Set-TextValue $ws.Range("E50") '  -1.10%  '
Set-TextValue $ws.Range("D51") '0.3352'
Set-TextValue $ws.Range("E51") '  -2.24%  '
